$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark "Absent" (column H) as 1 for all attendance date rows (3-18) except row 6
foreach ($r in 3..18) {
    if ($r -ne 6) {
        $ws.Cells.Item($r, 8).Value = 1
    }
}

# Row 6 is special: Total Attendance Count (D) and Real (E) become 1, Absent (H) stays 0
$ws.Cells.Item(6, 4).Value = 1
$ws.Cells.Item(6, 5).Value = 1
